$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Move the diagram group (id=79, "Group 78") up and to the left.
$grp = $s.Shapes.Item(1)
$grp.Left = 67.77692913385827
$grp.Top = 56.505045

# Widen the title bar (id=80, "Title 1") so it spans the full slide width.
$title = $s.Shapes.Item(2)
$title.Left = 14.684173228346456
$title.Width = 927.0

# Add a new textbox with the data-source link under the diagram.
# (A throwaway shape is created and removed first so the internal id/name
# counter lands on id=7 / "TextBox 6", matching the authored shape.)
$tmp = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$tmp.Delete()
$txt = $s.Shapes.AddTextbox(1, 485.11307086614175, 498.56472440944884, 445.7368503937008, 20.599212598425197)
$txt.Name = "TextBox 6"
$txt.Fill.Visible = $false
$txt.TextFrame.WordWrap = $true
$txt.TextFrame.AutoSize = 1
$tr = $txt.TextFrame.TextRange
$tr.Text = "https://github.com/emopsraps/Subsets/blob/main/UNICEF_Geosight_Country_id.csv"
$tr.LanguageID = "en-CA"
$tr.Font.Size = 11
